$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update database: swap the "Periodo Mora" (E) and "Valor Mora" (F) values
# between row 16 and row 18 (row 17 stays the same).
$ws.Range("E16").Value = "2107"
$ws.Range("F16").Value = 36341

$ws.Range("E18").Value = "2201"
$ws.Range("F18").Value = 30284
